# Updated yearly financial figures (Income Statement, Balance Sheet,
# and Cash Flow Statement) on sheet "EC" - "Doing Updates for Financials"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EC")

# Income Statement updates
# Row 8: D8=17667300, E8=15274400, F8=16669100, G8=21111000, H8=22537200, I8=21397500, J8=20069100
$ws.Cells.Item(8, 4).Value = 17667300
$ws.Cells.Item(8, 5).Value = 15274400
$ws.Cells.Item(8, 6).Value = 16669100
$ws.Cells.Item(8, 7).Value = 21111000
$ws.Cells.Item(8, 8).Value = 22537200
$ws.Cells.Item(8, 9).Value = 21397500
$ws.Cells.Item(8, 10).Value = 20069100

# Row 9: D9=11805900, E9=10955700, F9=11838200, G9=13947000, H9=26505000, I9=12273200, J9=10726200
$ws.Cells.Item(9, 4).Value = 11805900
$ws.Cells.Item(9, 5).Value = 10955700
$ws.Cells.Item(9, 6).Value = 11838200
$ws.Cells.Item(9, 7).Value = 13947000
$ws.Cells.Item(9, 8).Value = 26505000
$ws.Cells.Item(9, 9).Value = 12273200
$ws.Cells.Item(9, 10).Value = 10726200

# Row 10: D10=5861400, E10=4318600, F10=4830900, G10=7164000, H10=-3967800, I10=9124300, J10=9342800
$ws.Cells.Item(10, 4).Value = 5861400
$ws.Cells.Item(10, 5).Value = 4318600
$ws.Cells.Item(10, 6).Value = 4830900
$ws.Cells.Item(10, 7).Value = 7164000
$ws.Cells.Item(10, 8).Value = -3967800
$ws.Cells.Item(10, 9).Value = 9124300
$ws.Cells.Item(10, 10).Value = 9342800

# Row 12: D12=429400, E12=233100, F12=507000, G12=824400, H12=433600
$ws.Cells.Item(12, 4).Value = 429400
$ws.Cells.Item(12, 5).Value = 233100
$ws.Cells.Item(12, 6).Value = 507000
$ws.Cells.Item(12, 7).Value = 824400
$ws.Cells.Item(12, 8).Value = 433600

# Row 14: D14=-560500, E14=233300, F14=2584900, G14=672300
$ws.Cells.Item(14, 4).Value = -560500
$ws.Cells.Item(14, 5).Value = 233300
$ws.Cells.Item(14, 6).Value = 2584900
$ws.Cells.Item(14, 7).Value = 672300

# Row 15: D15=47800, E15=71400, F15=54600, G15=158100, H15=121600
$ws.Cells.Item(15, 4).Value = 47800
$ws.Cells.Item(15, 5).Value = 71400
$ws.Cells.Item(15, 6).Value = 54600
$ws.Cells.Item(15, 7).Value = 158100
$ws.Cells.Item(15, 8).Value = 121600

# Row 17: D17=12705800, E17=12633400, F17=16203200, G17=16487300, H17=15550100, I17=13859300, J17=12493400
$ws.Cells.Item(17, 4).Value = 12705800
$ws.Cells.Item(17, 5).Value = 12633400
$ws.Cells.Item(17, 6).Value = 16203200
$ws.Cells.Item(17, 7).Value = 16487300
$ws.Cells.Item(17, 8).Value = 15550100
$ws.Cells.Item(17, 9).Value = 13859300
$ws.Cells.Item(17, 10).Value = 12493400

# Row 18: D18=4961500, E18=2641000, F18=465900, G18=4623700, H18=6987100, I18=7538200, J18=7575600
$ws.Cells.Item(18, 4).Value = 4961500
$ws.Cells.Item(18, 5).Value = 2641000
$ws.Cells.Item(18, 6).Value = 465900
$ws.Cells.Item(18, 7).Value = 4623700
$ws.Cells.Item(18, 8).Value = 6987100
$ws.Cells.Item(18, 9).Value = 7538200
$ws.Cells.Item(18, 10).Value = 7575600

# Income Statement (continued) / Net Income reconciliation
# Row 20: D20=-26200, E20=503000, F20=-692400, G20=-824700, H20=204800, I20=-365900, J20=-69500
$ws.Cells.Item(20, 4).Value = -26200
$ws.Cells.Item(20, 5).Value = 503000
$ws.Cells.Item(20, 6).Value = -692400
$ws.Cells.Item(20, 7).Value = -824700
$ws.Cells.Item(20, 8).Value = 204800
$ws.Cells.Item(20, 9).Value = -365900
$ws.Cells.Item(20, 10).Value = -69500

# Row 21: D21=7580600, E21=5573400, F21=1940100, G21=7207700, H21=8957700, I21=8765000, J21=9141300
$ws.Cells.Item(21, 4).Value = 7580600
$ws.Cells.Item(21, 5).Value = 5573400
$ws.Cells.Item(21, 6).Value = 1940100
$ws.Cells.Item(21, 7).Value = 7207700
$ws.Cells.Item(21, 8).Value = 8957700
$ws.Cells.Item(21, 9).Value = 8765000
$ws.Cells.Item(21, 10).Value = 9141300

# Row 22: D22=763500, E22=884800, F22=566000, G22=241400, H22=189600
$ws.Cells.Item(22, 4).Value = 763500
$ws.Cells.Item(22, 5).Value = 884800
$ws.Cells.Item(22, 6).Value = 566000
$ws.Cells.Item(22, 7).Value = 241400
$ws.Cells.Item(22, 8).Value = 189600

# Row 23: D23=4171800, E23=2259100, F23=-792400, G23=3557700, H23=7002300, I23=7172300, J23=7506100
$ws.Cells.Item(23, 4).Value = 4171800
$ws.Cells.Item(23, 5).Value = 2259100
$ws.Cells.Item(23, 6).Value = -792400
$ws.Cells.Item(23, 7).Value = 3557700
$ws.Cells.Item(23, 8).Value = 7002300
$ws.Cells.Item(23, 9).Value = 7172300
$ws.Cells.Item(23, 10).Value = 7506100

# Row 24: D24=1803200, E24=1489800, F24=194100, G24=1526100, H24=2588400, I24=2408300, J24=2687700
$ws.Cells.Item(24, 4).Value = 1803200
$ws.Cells.Item(24, 5).Value = 1489800
$ws.Cells.Item(24, 6).Value = 194100
$ws.Cells.Item(24, 7).Value = 1526100
$ws.Cells.Item(24, 8).Value = 2588400
$ws.Cells.Item(24, 9).Value = 2408300
$ws.Cells.Item(24, 10).Value = 2687700

# Row 26: D26=2368600, E26=769400, F26=-986500, G26=2031500, H26=4413900, I26=4764000, J26=4818400
$ws.Cells.Item(26, 4).Value = 2368600
$ws.Cells.Item(26, 5).Value = 769400
$ws.Cells.Item(26, 6).Value = -986500
$ws.Cells.Item(26, 7).Value = 2031500
$ws.Cells.Item(26, 8).Value = 4413900
$ws.Cells.Item(26, 9).Value = 4764000
$ws.Cells.Item(26, 10).Value = 4818400

# Row 27: D27=2118500, E27=500700, F27=-1276100, G27=1832200, H27=4194100, I27=4702600, J27=4741500
$ws.Cells.Item(27, 4).Value = 2118500
$ws.Cells.Item(27, 5).Value = 500700
$ws.Cells.Item(27, 6).Value = -1276100
$ws.Cells.Item(27, 7).Value = 1832200
$ws.Cells.Item(27, 8).Value = 4194100
$ws.Cells.Item(27, 9).Value = 4702600
$ws.Cells.Item(27, 10).Value = 4741500

# Row 32: D32=26200, E32=-503000, F32=692400, G32=824700, H32=-204800, I32=365900, J32=69500
$ws.Cells.Item(32, 4).Value = 26200
$ws.Cells.Item(32, 5).Value = -503000
$ws.Cells.Item(32, 6).Value = 692400
$ws.Cells.Item(32, 7).Value = 824700
$ws.Cells.Item(32, 8).Value = -204800
$ws.Cells.Item(32, 9).Value = 365900
$ws.Cells.Item(32, 10).Value = 69500

# Row 33: D33=2118500, E33=500700, F33=-1276100, G33=1832200, H33=4194100, I33=4702600, J33=4741500
$ws.Cells.Item(33, 4).Value = 2118500
$ws.Cells.Item(33, 5).Value = 500700
$ws.Cells.Item(33, 6).Value = -1276100
$ws.Cells.Item(33, 7).Value = 1832200
$ws.Cells.Item(33, 8).Value = 4194100
$ws.Cells.Item(33, 9).Value = 4702600
$ws.Cells.Item(33, 10).Value = 4741500

# Row 35: D35=2118500, E35=500700, F35=-1276100, G35=1832200, H35=4194100, I35=4702600, J35=4741500
$ws.Cells.Item(35, 4).Value = 2118500
$ws.Cells.Item(35, 5).Value = 500700
$ws.Cells.Item(35, 6).Value = -1276100
$ws.Cells.Item(35, 7).Value = 1832200
$ws.Cells.Item(35, 8).Value = 4194100
$ws.Cells.Item(35, 9).Value = 4702600
$ws.Cells.Item(35, 10).Value = 4741500

# Balance Sheet - Assets
# Row 41: D41=1755700, E41=1062500, F41=1435100, G41=4208400, H41=5339300, I41=2551100, J41=2263500
$ws.Cells.Item(41, 4).Value = 1755700
$ws.Cells.Item(41, 5).Value = 1062500
$ws.Cells.Item(41, 6).Value = 1435100
$ws.Cells.Item(41, 7).Value = 4208400
$ws.Cells.Item(41, 8).Value = 5339300
$ws.Cells.Item(41, 9).Value = 2551100
$ws.Cells.Item(41, 10).Value = 2263500

# Row 42: D42=1736700, E42=3346300, F42=1058700, G42=1915600, H42=1344500, I42=855500, J42=297000
$ws.Cells.Item(42, 4).Value = 1736700
$ws.Cells.Item(42, 5).Value = 3346300
$ws.Cells.Item(42, 6).Value = 1058700
$ws.Cells.Item(42, 7).Value = 1915600
$ws.Cells.Item(42, 8).Value = 1344500
$ws.Cells.Item(42, 9).Value = 855500
$ws.Cells.Item(42, 10).Value = 297000

# Row 43: D43=2341000, E43=1946900, F43=2684600, G43=3692700, H43=1976600, I43=1944300, J43=1752000
$ws.Cells.Item(43, 4).Value = 2341000
$ws.Cells.Item(43, 5).Value = 1946900
$ws.Cells.Item(43, 6).Value = 2684600
$ws.Cells.Item(43, 7).Value = 3692700
$ws.Cells.Item(43, 8).Value = 1976600
$ws.Cells.Item(43, 9).Value = 1944300
$ws.Cells.Item(43, 10).Value = 1752000

# Row 44: D44=1472400, E44=1229400, F44=978500, G44=2153900, H44=2187300, I44=854500, J44=862800
$ws.Cells.Item(44, 4).Value = 1472400
$ws.Cells.Item(44, 5).Value = 1229400
$ws.Cells.Item(44, 6).Value = 978500
$ws.Cells.Item(44, 7).Value = 2153900
$ws.Cells.Item(44, 8).Value = 2187300
$ws.Cells.Item(44, 9).Value = 854500
$ws.Cells.Item(44, 10).Value = 862800

# Row 45: D45=125800, E45=136200, F45=279300, G45=1775200, H45=3042100, I45=789300, J45=523800
$ws.Cells.Item(45, 4).Value = 125800
$ws.Cells.Item(45, 5).Value = 136200
$ws.Cells.Item(45, 6).Value = 279300
$ws.Cells.Item(45, 7).Value = 1775200
$ws.Cells.Item(45, 8).Value = 3042100
$ws.Cells.Item(45, 9).Value = 789300
$ws.Cells.Item(45, 10).Value = 523800

# Row 46: D46=7431700, E46=7721300, F46=6436300, G46=6605000, H46=9381600, I46=6452100, J46=5699000
$ws.Cells.Item(46, 4).Value = 7431700
$ws.Cells.Item(46, 5).Value = 7721300
$ws.Cells.Item(46, 6).Value = 6436300
$ws.Cells.Item(46, 7).Value = 6605000
$ws.Cells.Item(46, 8).Value = 9381600
$ws.Cells.Item(46, 9).Value = 6452100
$ws.Cells.Item(46, 10).Value = 5699000

# Row 47: D47=1815500, E47=1169100, F47=1207300, G47=1418800, H47=1693800, I47=4291400, J47=2246100
$ws.Cells.Item(47, 4).Value = 1815500
$ws.Cells.Item(47, 5).Value = 1169100
$ws.Cells.Item(47, 6).Value = 1207300
$ws.Cells.Item(47, 7).Value = 1418800
$ws.Cells.Item(47, 8).Value = 1693800
$ws.Cells.Item(47, 9).Value = 4291400
$ws.Cells.Item(47, 10).Value = 2246100

# Row 48: D48=26453800, E48=27075200, F48=28503700, G48=52342400, H48=41267800, I48=26792300, J48=13214700
$ws.Cells.Item(48, 4).Value = 26453800
$ws.Cells.Item(48, 5).Value = 27075200
$ws.Cells.Item(48, 6).Value = 28503700
$ws.Cells.Item(48, 7).Value = 52342400
$ws.Cells.Item(48, 8).Value = 41267800
$ws.Cells.Item(48, 9).Value = 26792300
$ws.Cells.Item(48, 10).Value = 13214700

# Row 49: D49=415900, E49=381300, F49=418400, G49=612400, H49=1140300, I49=404600, J49=444300
$ws.Cells.Item(49, 4).Value = 415900
$ws.Cells.Item(49, 5).Value = 381300
$ws.Cells.Item(49, 6).Value = 418400
$ws.Cells.Item(49, 7).Value = 612400
$ws.Cells.Item(49, 8).Value = 1140300
$ws.Cells.Item(49, 9).Value = 404600
$ws.Cells.Item(49, 10).Value = 444300

# Row 52: D52=2499000, E52=2471400, F52=2793100, G52=13046000, H52=10772200, I52=1046400, J52=1086700
$ws.Cells.Item(52, 4).Value = 2499000
$ws.Cells.Item(52, 5).Value = 2471400
$ws.Cells.Item(52, 6).Value = 2793100
$ws.Cells.Item(52, 7).Value = 13046000
$ws.Cells.Item(52, 8).Value = 10772200
$ws.Cells.Item(52, 9).Value = 1046400
$ws.Cells.Item(52, 10).Value = 1086700

# Row 54: D54=38615900, E54=38818400, F54=39358700, G54=35468100, H54=42377000, I54=26086200, J54=22690900
$ws.Cells.Item(54, 4).Value = 38615900
$ws.Cells.Item(54, 5).Value = 38818400
$ws.Cells.Item(54, 6).Value = 39358700
$ws.Cells.Item(54, 7).Value = 35468100
$ws.Cells.Item(54, 8).Value = 42377000
$ws.Cells.Item(54, 9).Value = 26086200
$ws.Cells.Item(54, 10).Value = 22690900

# Balance Sheet - Liabilities
# Row 57: D57=1628500, E57=1494300, F57=1593600, G57=4792700, H57=4352500, I57=2314700, J57=1613900
$ws.Cells.Item(57, 4).Value = 1628500
$ws.Cells.Item(57, 5).Value = 1494300
$ws.Cells.Item(57, 6).Value = 1593600
$ws.Cells.Item(57, 7).Value = 4792700
$ws.Cells.Item(57, 8).Value = 4352500
$ws.Cells.Item(57, 9).Value = 2314700
$ws.Cells.Item(57, 10).Value = 1613900

# Row 58: D58=1646200, E58=1320400, F58=1463600, G58=1125600, H58=288400, I58=690900, J58=302800
$ws.Cells.Item(58, 4).Value = 1646200
$ws.Cells.Item(58, 5).Value = 1320400
$ws.Cells.Item(58, 6).Value = 1463600
$ws.Cells.Item(58, 7).Value = 1125600
$ws.Cells.Item(58, 8).Value = 288400
$ws.Cells.Item(58, 9).Value = 690900
$ws.Cells.Item(58, 10).Value = 302800

# Row 59: D59=2116200, E59=2429300, F59=2524800, G59=2514300, H59=5128400, I59=3115500, J59=2558000
$ws.Cells.Item(59, 4).Value = 2116200
$ws.Cells.Item(59, 5).Value = 2429300
$ws.Cells.Item(59, 6).Value = 2524800
$ws.Cells.Item(59, 7).Value = 2514300
$ws.Cells.Item(59, 8).Value = 5128400
$ws.Cells.Item(59, 9).Value = 3115500
$ws.Cells.Item(59, 10).Value = 2558000

# Row 60: D60=5390900, E60=5244000, F60=5582000, G60=5390000, H60=7115400, I60=5693400, J60=4474700
$ws.Cells.Item(60, 4).Value = 5390900
$ws.Cells.Item(60, 5).Value = 5244000
$ws.Cells.Item(60, 6).Value = 5582000
$ws.Cells.Item(60, 7).Value = 5390000
$ws.Cells.Item(60, 8).Value = 7115400
$ws.Cells.Item(60, 9).Value = 5693400
$ws.Cells.Item(60, 10).Value = 4474700

# Row 61: D61=12289100, E61=15390700, F61=15567900, G61=10087700, H61=6855700, I61=4280700, J61=2799300
$ws.Cells.Item(61, 4).Value = 12289100
$ws.Cells.Item(61, 5).Value = 15390700
$ws.Cells.Item(61, 6).Value = 15567900
$ws.Cells.Item(61, 7).Value = 10087700
$ws.Cells.Item(61, 8).Value = 6855700
$ws.Cells.Item(61, 9).Value = 4280700
$ws.Cells.Item(61, 10).Value = 2799300

# Row 62: D62=5005900, E62=3681400, F62=3734700, G62=7237200, H62=4281900, I62=3458400, J62=3147000
$ws.Cells.Item(62, 4).Value = 5005900
$ws.Cells.Item(62, 5).Value = 3681400
$ws.Cells.Item(62, 6).Value = 3734700
$ws.Cells.Item(62, 7).Value = 7237200
$ws.Cells.Item(62, 8).Value = 4281900
$ws.Cells.Item(62, 9).Value = 3458400
$ws.Cells.Item(62, 10).Value = 3147000

# Row 66: D66=23288300, E66=24842700, F66=25484600, G66=20101300, H66=19618800, I66=14038700, J66=11153200
$ws.Cells.Item(66, 4).Value = 23288300
$ws.Cells.Item(66, 5).Value = 24842700
$ws.Cells.Item(66, 6).Value = 25484600
$ws.Cells.Item(66, 7).Value = 20101300
$ws.Cells.Item(66, 8).Value = 19618800
$ws.Cells.Item(66, 9).Value = 14038700
$ws.Cells.Item(66, 10).Value = 11153200

# Balance Sheet - Stockholders' Equity
# Row 72: D72=3163800, E72=1255700, F72=449200, G72=27931600, H72=25460100, I72=16053700, J72=7271700
$ws.Cells.Item(72, 4).Value = 3163800
$ws.Cells.Item(72, 5).Value = 1255700
$ws.Cells.Item(72, 6).Value = 449200
$ws.Cells.Item(72, 7).Value = 27931600
$ws.Cells.Item(72, 8).Value = 25460100
$ws.Cells.Item(72, 9).Value = 16053700
$ws.Cells.Item(72, 10).Value = 7271700

# Row 76: D76=15327600, E76=13975600, F76=13874100, G76=15366800, H76=22758100, I76=12047500, J76=11537700
$ws.Cells.Item(76, 4).Value = 15327600
$ws.Cells.Item(76, 5).Value = 13975600
$ws.Cells.Item(76, 6).Value = 13874100
$ws.Cells.Item(76, 7).Value = 15366800
$ws.Cells.Item(76, 8).Value = 22758100
$ws.Cells.Item(76, 9).Value = 12047500
$ws.Cells.Item(76, 10).Value = 11537700

# Cash Flow Statement - Operating Activities
# Row 81: D81=2118500, E81=500700, F81=-1276100, G81=1832200, H81=4194100, I81=4702600, J81=4741500
$ws.Cells.Item(81, 4).Value = 2118500
$ws.Cells.Item(81, 5).Value = 500700
$ws.Cells.Item(81, 6).Value = -1276100
$ws.Cells.Item(81, 7).Value = 1832200
$ws.Cells.Item(81, 8).Value = 4194100
$ws.Cells.Item(81, 9).Value = 4702600
$ws.Cells.Item(81, 10).Value = 4741500

# Row 83: D83=2645300, E83=2429500, F83=2166500, G83=3408700, H83=1765800, I83=1592600, J83=1635100
$ws.Cells.Item(83, 4).Value = 2645300
$ws.Cells.Item(83, 5).Value = 2429500
$ws.Cells.Item(83, 6).Value = 2166500
$ws.Cells.Item(83, 7).Value = 3408700
$ws.Cells.Item(83, 8).Value = 1765800
$ws.Cells.Item(83, 9).Value = 1592600
$ws.Cells.Item(83, 10).Value = 1635100

# Row 89: D89=5431600, E89=4620200, F89=3331700, G89=5303400, H89=5607300, I89=6440700, J89=7091300
$ws.Cells.Item(89, 4).Value = 5431600
$ws.Cells.Item(89, 5).Value = 4620200
$ws.Cells.Item(89, 6).Value = 3331700
$ws.Cells.Item(89, 7).Value = 5303400
$ws.Cells.Item(89, 8).Value = 5607300
$ws.Cells.Item(89, 9).Value = 6440700
$ws.Cells.Item(89, 10).Value = 7091300

# Cash Flow Statement - Investing/Financing Activities
# Row 91: D91=-1852700, E91=-1845800, F91=-4524600, G91=-4467800, H91=-2547200, I91=-2992000, J91=-3232100
$ws.Cells.Item(91, 4).Value = -1852700
$ws.Cells.Item(91, 5).Value = -1845800
$ws.Cells.Item(91, 6).Value = -4524600
$ws.Cells.Item(91, 7).Value = -4467800
$ws.Cells.Item(91, 8).Value = -2547200
$ws.Cells.Item(91, 9).Value = -2992000
$ws.Cells.Item(91, 10).Value = -3232100

# Row 94: D94=-1402800, E94=-3068400, F94=-3694400, G94=-3702000, H94=-3092000, I94=-5124200, J94=-4673200
$ws.Cells.Item(94, 4).Value = -1402800
$ws.Cells.Item(94, 5).Value = -3068400
$ws.Cells.Item(94, 6).Value = -3694400
$ws.Cells.Item(94, 7).Value = -3702000
$ws.Cells.Item(94, 8).Value = -3092000
$ws.Cells.Item(94, 9).Value = -5124200
$ws.Cells.Item(94, 10).Value = -4673200

# Row 96: D96=-481500, E96=-547900, F96=-1757900, G96=-572600, H96=-4662500, I96=-2694200, J96=-1890200
$ws.Cells.Item(96, 4).Value = -481500
$ws.Cells.Item(96, 5).Value = -547900
$ws.Cells.Item(96, 6).Value = -1757900
$ws.Cells.Item(96, 7).Value = -572600
$ws.Cells.Item(96, 8).Value = -4662500
$ws.Cells.Item(96, 9).Value = -2694200
$ws.Cells.Item(96, 10).Value = -1890200

# Cash Flow Statement - Effect of exchange rate / change in cash
# Row 100: D100=-4084500, E100=-884200, F100=-445500, G100=-2350800, H100=-2264900, I100=-1106300, J100=-1400900
$ws.Cells.Item(100, 4).Value = -4084500
$ws.Cells.Item(100, 5).Value = -884200
$ws.Cells.Item(100, 6).Value = -445500
$ws.Cells.Item(100, 7).Value = -2350800
$ws.Cells.Item(100, 8).Value = -2264900
$ws.Cells.Item(100, 9).Value = -1106300
$ws.Cells.Item(100, 10).Value = -1400900

# Row 101: D101=-92900, E101=-72400, F101=466600, G101=369700, H101=37800, I101=77400, J101=-5100
$ws.Cells.Item(101, 4).Value = -92900
$ws.Cells.Item(101, 5).Value = -72400
$ws.Cells.Item(101, 6).Value = 466600
$ws.Cells.Item(101, 7).Value = 369700
$ws.Cells.Item(101, 8).Value = 37800
$ws.Cells.Item(101, 9).Value = 77400
$ws.Cells.Item(101, 10).Value = -5100

# Row 102: D102=-148700, E102=595200, F102=-341700, G102=-379800, H102=288200, I102=287600, J102=1012100
$ws.Cells.Item(102, 4).Value = -148700
$ws.Cells.Item(102, 5).Value = 595200
$ws.Cells.Item(102, 6).Value = -341700
$ws.Cells.Item(102, 7).Value = -379800
$ws.Cells.Item(102, 8).Value = 288200
$ws.Cells.Item(102, 9).Value = 287600
$ws.Cells.Item(102, 10).Value = 1012100
